$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 header values
$ws.Range("B1").Value = 16
$ws.Range("C1").Value = 20
$ws.Range("D1").Value = 16
$ws.Range("E1").Value = 20

# Row 2 (CON) values - B2 and D2 are deleted (cleared), C2 and E2 updated
$ws.Range("B2").ClearContents()
$ws.Range("C2").Value = -5.4409116359273781
$ws.Range("D2").ClearContents()
$ws.Range("E2").Value = -3.8984633129948634

# Row 3 (STR) values
$ws.Range("B3").Value = -6.2925999639750998
$ws.Range("C3").Value = 1.1129731498115163
$ws.Range("D3").Value = -5.548955410987837
$ws.Range("E3").Value = 9.7698892100261858

# Update selection to match the new reduced selection range
$ws.Range("B1:E3").Select()
